# Fixed one-based index with Excel parser.
# Replace the month names in column B with person names (rows 2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$names = @("Greg", "Don", "Hilda", "Fran", "Eddie", "Paul", "Peter", "William", "Oliver", "Owen")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Update the active cell selection to B12 (one row past the last data row).
$ws.Range("B12").Select()
